{"js": "// The author's edit (per the diff) leaves the visible text of the document\n// unchanged -- \"... send the HTTPS clone URL to the TA in an email.\" still\n// reads the same.  What actually moved is Word's hidden \"_GoBack\" bookmark\n// (the marker Word drops at the site of the most recent edit, and that\n// \"Resume Reading\"/Shift+F5 jumps back to).  It used to sit right after\n// \"Due: Tomorrow at 5pm\"; the diff shows it removed from there and\n// re-inserted immediately before the final period of the sentence that\n// ends \"... to the TA in an email.\", splitting that run's text node into\n// \"... an email\" + \".\" around the relocated bookmark.\n//\n// Reproduce that: drop the old \"_GoBack\" bookmark, find the tail of that\n// sentence, and plant a fresh \"_GoBack\" bookmark collapsed right before the\n// trailing period.\n\n// 1. Remove the existing \"_GoBack\" bookmark (after \"Due: Tomorrow at 5pm\").\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2. Locate the sentence fragment immediately preceding the final period,\n//    so a collapsed range at its end sits exactly between \"email\" and \".\".\nconst results = context.document.body.search(\"TA in an email\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const found = results.items[0];\n  const insertionPoint = found.getRange(\"End\");\n\n  // 3. Re-create \"_GoBack\" as a collapsed bookmark at that exact spot.\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The author's edit (per the diff) leaves the visible text of the document\n# unchanged -- \"... send the HTTPS clone URL to the TA in an email.\" still\n# reads the same.  What actually moved is Word's hidden \"_GoBack\" bookmark\n# (the marker Word drops at the site of the most recent edit, and that\n# \"Resume Reading\"/Shift+F5 jumps back to).  It used to sit right after\n# \"Due: Tomorrow at 5pm\"; the diff shows it removed from there and\n# re-inserted immediately before the final period of the sentence that\n# ends \"... to the TA in an email.\", splitting that run's text node into\n# \"... an email\" + \".\" around the relocated bookmark.\n#\n# Reproduce that: drop the old \"_GoBack\" bookmark, find the tail of that\n# sentence, and plant a fresh \"_GoBack\" bookmark collapsed right before the\n# trailing period.\n\n$d = $word.ActiveDocument\n\n# 1. Remove the existing \"_GoBack\" bookmark (after \"Due: Tomorrow at 5pm\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Locate the sentence fragment immediately preceding the final period,\n#    so a collapsed range at its end sits exactly between \"email\" and \".\".\n$find = $d.Content\n$find.Find.Text = \"TA in an email\"\n$found = $find.Find.Execute()\n\nif ($found) {\n    # Collapse the found range to its end point (right before the period).\n    $insertionPoint = $d.Range($find.End, $find.End)\n\n    # 3. Re-create \"_GoBack\" as a collapsed bookmark at that exact spot.\n    $d.Bookmarks.Add(\"_GoBack\", $insertionPoint)\n}\n"}
